# CDS Input file updates
# - Replace the "Participant ID" Cypher query text (column B, row 2 / ParticipantsTab)
#   with an updated version that also considers diagnosis + genomic_info and
#   returns samples sorted via apoc.coll.sort.
# - Row 2 grows taller to fit the longer query text.
# - Selection/viewport moves to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newParticipantQuery = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE f.file_type in ['HTML']
with p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id LIMIT 100
"@

# Update the query text held in B2 (ParticipantsTab row).
$ws.Range("B2").Value = $newParticipantQuery

# The row needs to grow to accommodate the extra lines of text.
$ws.Rows.Item(2).RowHeight = 279

# Move the visible selection to B4, as recorded in the saved view state.
$ws.Range("B4").Select() | Out-Null
